# Update "想去人数" (F column) values on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to reflect newly generated counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 292
$ws1.Range("F6").Value = 1113
$ws1.Range("F9").Value = 113
$ws1.Range("F14").Value = 448
$ws1.Range("F15").Value = 1388
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 118
$ws1.Range("F20").Value = 69
$ws1.Range("F22").Value = 1010
$ws1.Range("F24").Value = 241
$ws1.Range("F25").Value = 26
$ws1.Range("F26").Value = 5966
$ws1.Range("F28").Value = 124
$ws1.Range("F29").Value = 118
$ws1.Range("F31").Value = 14660
$ws1.Range("F32").Value = 1455
$ws1.Range("F33").Value = 226
$ws1.Range("F34").Value = 107
$ws1.Range("F36").Value = 9484
$ws1.Range("F37").Value = 641
$ws1.Range("F39").Value = 159
$ws1.Range("F40").Value = 364
$ws1.Range("F41").Value = 112

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 292
$ws4.Range("F6").Value = 1113
$ws4.Range("F9").Value = 113
$ws4.Range("F14").Value = 448
$ws4.Range("F15").Value = 1388
$ws4.Range("F16").Value = 124
$ws4.Range("F17").Value = 118
$ws4.Range("F21").Value = 69
$ws4.Range("F24").Value = 1010
$ws4.Range("F26").Value = 241
$ws4.Range("F27").Value = 26
$ws4.Range("F29").Value = 5966
$ws4.Range("F31").Value = 124
$ws4.Range("F32").Value = 118
$ws4.Range("F34").Value = 14660
$ws4.Range("F35").Value = 1455
$ws4.Range("F36").Value = 226
$ws4.Range("F37").Value = 107
$ws4.Range("F39").Value = 9484
$ws4.Range("F40").Value = 641
$ws4.Range("F42").Value = 159
$ws4.Range("F43").Value = 364
$ws4.Range("F44").Value = 112
